$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.039
$ws.Range("E2").Value = 0.07199999999999999
$ws.Range("F2").Value = 0.046
$ws.Range("D3").Value = 0.344
$ws.Range("E3").Value = 0.389
$ws.Range("F3").Value = 0.355
$ws.Range("D4").Value = 0.051
$ws.Range("E4").Value = 0.056
$ws.Range("F4").Value = 0.052
$ws.Range("D5").Value = 4.764
$ws.Range("E5").Value = 7.18
$ws.Range("F5").Value = 5.225
$ws.Range("D6").Value = 0.473
$ws.Range("E6").Value = 0.509
$ws.Range("F6").Value = 0.484
$ws.Range("H6").Value = 5005
$ws.Range("D7").Value = 2.093
$ws.Range("E7").Value = 2.191
$ws.Range("F7").Value = 2.115
$ws.Range("D8").Value = 6.275
$ws.Range("E8").Value = 6.648
$ws.Range("F8").Value = 6.347
$ws.Range("D9").Value = 3.855
$ws.Range("E9").Value = 3.939
$ws.Range("F9").Value = 3.859
$ws.Range("H9").Value = 11648
$ws.Range("D10").Value = 4.005
$ws.Range("E10").Value = 4.046
$ws.Range("F10").Value = 3.983
$ws.Range("H10").Value = 11648
$ws.Range("D11").Value = 6.763
$ws.Range("E11").Value = 7.447
$ws.Range("F11").Value = 6.781
$ws.Range("H11").Value = 11648
$ws.Range("D12").Value = 3.819
$ws.Range("E12").Value = 3.86
$ws.Range("F12").Value = 3.784
$ws.Range("D13").Value = 3.573
$ws.Range("E13").Value = 3.843
$ws.Range("F13").Value = 3.522
$ws.Range("D14").Value = 2.851
$ws.Range("E14").Value = 2.958
$ws.Range("F14").Value = 2.858
$ws.Range("D15").Value = 8.618
$ws.Range("E15").Value = 9.237
$ws.Range("F15").Value = 8.612
$ws.Range("D16").Value = 8.484
$ws.Range("E16").Value = 9.004
$ws.Range("F16").Value = 8.491
$ws.Range("D17").Value = 4.329
$ws.Range("E17").Value = 4.649
$ws.Range("F17").Value = 4.335
$ws.Range("D18").Value = 6.468
$ws.Range("E18").Value = 7.074
$ws.Range("F18").Value = 6.612
$ws.Range("D19").Value = 4.01
$ws.Range("E19").Value = 4.032
$ws.Range("F19").Value = 3.958
$ws.Range("H19").Value = 11648
$ws.Range("D20").Value = 4.311
$ws.Range("E20").Value = 4.557
$ws.Range("F20").Value = 4.301
$ws.Range("H20").Value = 11648
$ws.Range("D21").Value = 6.957
$ws.Range("E21").Value = 7.808
$ws.Range("F21").Value = 7.014
$ws.Range("H21").Value = 11648
$ws.Range("D22").Value = 0.049
$ws.Range("E22").Value = 0.055
$ws.Range("F22").Value = 0.05
$ws.Range("D23").Value = 0.313
$ws.Range("E23").Value = 0.347
$ws.Range("F23").Value = 0.322
$ws.Range("H23").Value = 2688
$ws.Range("D24").Value = 0.042
$ws.Range("E24").Value = 0.044
$ws.Range("F24").Value = 0.043
$ws.Range("H24").Value = 56
$ws.Range("D25").Value = 6.056
$ws.Range("E25").Value = 6.72
$ws.Range("F25").Value = 6.198
$ws.Range("D26").Value = 0.333
$ws.Range("E26").Value = 0.36
$ws.Range("F26").Value = 0.34
$ws.Range("H26").Value = 3080
$ws.Range("D27").Value = 3.309
$ws.Range("E27").Value = 3.405
$ws.Range("F27").Value = 3.334
$ws.Range("H27").Value = 3717
$ws.Range("D28").Value = 3.126
$ws.Range("E28").Value = 3.528
$ws.Range("F28").Value = 3.209
$ws.Range("D29").Value = 0.24
$ws.Range("E29").Value = 0.263
$ws.Range("F29").Value = 0.245
$ws.Range("H29").Value = 2695
$ws.Range("E30").Value = 0.01
$ws.Range("D31").Value = 1.125
$ws.Range("E31").Value = 1.155
$ws.Range("F31").Value = 1.119
$ws.Range("D32").Value = 1.264
$ws.Range("E32").Value = 1.298
$ws.Range("F32").Value = 1.252
$ws.Range("D33").Value = 2.339
$ws.Range("E33").Value = 2.633
$ws.Range("F33").Value = 2.363
$ws.Range("D34").Value = 1.311
$ws.Range("E34").Value = 1.373
$ws.Range("F34").Value = 1.319
$ws.Range("D35").Value = 6.312
$ws.Range("E35").Value = 6.44
$ws.Range("F35").Value = 6.275
$ws.Range("D36").Value = 1.247
$ws.Range("E36").Value = 1.309
$ws.Range("F36").Value = 1.246
$ws.Range("D37").Value = 2.939
$ws.Range("E37").Value = 3.006
$ws.Range("F37").Value = 2.942
$ws.Range("D38").Value = 3.369
$ws.Range("E38").Value = 5.713
$ws.Range("F38").Value = 3.905
$ws.Range("D39").Value = 1.929
$ws.Range("E39").Value = 2.269
$ws.Range("F39").Value = 2.027
$ws.Range("D40").Value = 2.354
$ws.Range("E40").Value = 2.413
$ws.Range("F40").Value = 2.366
$ws.Range("D41").Value = 2.377
$ws.Range("E41").Value = 2.405
$ws.Range("F41").Value = 2.376
